# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Header row (row 1): AD1=Wins, AE1=Losses, AF1=Ties --------------------
$headers = @("Wins", "Losses", "Ties")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 30 + $i   # AD=30, AE=31, AF=32
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]

    # Match the formatting used by the other header cells (bold, centered,
    # thin border) such as A1:AC1.
    $cell.Font.Bold = $true
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# -- Data rows (rows 2-47): season record for every player -----------------
$wins = 78
$losses = 84
$ties = 0

for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
